# Auto-generated edit script: refresh currentAveragePrice / Leve price & profit columns
# (H..N) for the rows whose source market data changed, per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 838.381
$ws.Range("I107").Value = 572.55554
$ws.Range("K107").Value = 572.55554
$ws.Range("M107").Value = 1347.44446
# Row 137
$ws.Range("H137").Value = 627922.25
$ws.Range("I137").Value = 3155.6667
$ws.Range("J137").Value = 2502222
$ws.Range("K137").Value = 9467.000100000001
$ws.Range("L137").Value = 7506666
$ws.Range("M137").Value = -6917.000100000001
$ws.Range("N137").Value = -7511766
# Row 138
$ws.Range("H138").Value = 5588.0557
$ws.Range("I138").Value = 1589
$ws.Range("J138").Value = 6730.643
$ws.Range("K138").Value = 4767
$ws.Range("L138").Value = 20191.929
$ws.Range("M138").Value = 373
$ws.Range("N138").Value = -30471.929

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17264.754
$ws.Range("I32").Value = 17976.895
$ws.Range("J32").Value = 14874
$ws.Range("K32").Value = 17976.895
$ws.Range("L32").Value = 14874
$ws.Range("M32").Value = -17689.895
$ws.Range("N32").Value = -15448
# Row 74
$ws.Range("H74").Value = 8622461
$ws.Range("I74").Value = 13158961
$ws.Range("J74").Value = 3111.2
$ws.Range("K74").Value = 13158961
$ws.Range("L74").Value = 3111.2
$ws.Range("M74").Value = -13158087
$ws.Range("N74").Value = -4859.2
# Row 77
$ws.Range("H77").Value = 8622461
$ws.Range("I77").Value = 13158961
$ws.Range("J77").Value = 3111.2
$ws.Range("K77").Value = 65794805
$ws.Range("L77").Value = 15556
$ws.Range("M77").Value = -65790437
$ws.Range("N77").Value = -24292

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 10000
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
# Row 141
$ws.Range("H141").Value = 110942
$ws.Range("J141").Value = 111107
$ws.Range("L141").Value = 111107
$ws.Range("N141").Value = -121467

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 15627971
$ws.Range("I31").Value = 21278008
$ws.Range("J31").Value = 7280
$ws.Range("K31").Value = 21278008
$ws.Range("L31").Value = 7280
$ws.Range("M31").Value = -21277713
$ws.Range("N31").Value = -7870
# Row 34
$ws.Range("H34").Value = 15627971
$ws.Range("I34").Value = 21278008
$ws.Range("J34").Value = 7280
$ws.Range("K34").Value = 21278008
$ws.Range("L34").Value = 7280
$ws.Range("M34").Value = -21277806
$ws.Range("N34").Value = -7684
# Row 132
$ws.Range("H132").Value = 19327410
$ws.Range("I132").Value = 23394194
$ws.Range("J132").Value = 10187.167
$ws.Range("K132").Value = 70182582
$ws.Range("L132").Value = 30561.501
$ws.Range("M132").Value = -70180052
$ws.Range("N132").Value = -35621.501

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 35863264
$ws.Range("I4").Value = 34968804
$ws.Range("J4").Value = 37858600
$ws.Range("K4").Value = 104906412
$ws.Range("L4").Value = 113575800
$ws.Range("M4").Value = -104906300
$ws.Range("N4").Value = -113576024
# Row 57
$ws.Range("H57").Value = 2233.3333
$ws.Range("I57").Value = 2233.3333
$ws.Range("K57").Value = 6699.999899999999
$ws.Range("M57").Value = -6140.999899999999
# Row 62
$ws.Range("H62").Value = 3999.5
$ws.Range("J62").Value = 3999.5
$ws.Range("L62").Value = 11998.5
$ws.Range("N62").Value = -13370.5
# Row 63
$ws.Range("H63").Value = 6377.3335
$ws.Range("I63").Value = 4171.6665
$ws.Range("J63").Value = 8583
$ws.Range("K63").Value = 12514.9995
$ws.Range("L63").Value = 25749
$ws.Range("M63").Value = -11765.9995
$ws.Range("N63").Value = -27247
# Row 64
$ws.Range("H64").Value = 18000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 18000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 54000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -54540
# Row 65
$ws.Range("H65").Value = 3999.5
$ws.Range("J65").Value = 3999.5
$ws.Range("L65").Value = 35995.5
$ws.Range("N65").Value = -42859.5
# Row 66
$ws.Range("H66").Value = 6377.3335
$ws.Range("I66").Value = 4171.6665
$ws.Range("J66").Value = 8583
$ws.Range("K66").Value = 37544.9985
$ws.Range("L66").Value = 77247
$ws.Range("M66").Value = -33800.9985
$ws.Range("N66").Value = -84735
# Row 67
$ws.Range("H67").Value = 18000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 18000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 54000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -55872
# Row 68
$ws.Range("H68").Value = 1828.5714
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1966.6666
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 5899.9998
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -7521.9998
# Row 70
$ws.Range("H70").Value = 5854.5
$ws.Range("I70").Value = 7563.5
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 22690.5
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -22375.5
$ws.Range("N70").Value = -15630
# Row 71
$ws.Range("H71").Value = 1828.5714
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1966.6666
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 17699.9994
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -25811.9994
# Row 73
$ws.Range("H73").Value = 5854.5
$ws.Range("I73").Value = 7563.5
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 22690.5
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -21598.5
$ws.Range("N73").Value = -17184
# Row 74
$ws.Range("H74").Value = 5833
$ws.Range("J74").Value = 6249.5
$ws.Range("L74").Value = 18748.5
$ws.Range("N74").Value = -20870.5
# Row 75
$ws.Range("H75").Value = 3789.7144
$ws.Range("J75").Value = 3789.7144
$ws.Range("L75").Value = 11369.1432
$ws.Range("N75").Value = -13365.1432
# Row 76
$ws.Range("H76").Value = 7500
$ws.Range("J76").Value = 7500
$ws.Range("L76").Value = 22500
$ws.Range("N76").Value = -23266
# Row 77
$ws.Range("H77").Value = 5833
$ws.Range("J77").Value = 6249.5
$ws.Range("L77").Value = 56245.5
$ws.Range("N77").Value = -66853.5
# Row 78
$ws.Range("H78").Value = 3789.7144
$ws.Range("J78").Value = 3789.7144
$ws.Range("L78").Value = 34107.4296
$ws.Range("N78").Value = -44091.4296
# Row 79
$ws.Range("H79").Value = 7500
$ws.Range("J79").Value = 7500
$ws.Range("L79").Value = 22500
$ws.Range("N79").Value = -25152
# Row 86
$ws.Range("H86").Value = 249
$ws.Range("I86").Value = 249
$ws.Range("K86").Value = 747
$ws.Range("M86").Value = 439
# Row 89
$ws.Range("H89").Value = 249
$ws.Range("I89").Value = 249
$ws.Range("K89").Value = 2241
$ws.Range("M89").Value = 3687
# Row 107
$ws.Range("H107").Value = 680.0833
$ws.Range("I107").Value = 335
$ws.Range("J107").Value = 852.625
$ws.Range("K107").Value = 1005
$ws.Range("L107").Value = 2557.875
$ws.Range("M107").Value = 915
$ws.Range("N107").Value = -6397.875
# Row 115
$ws.Range("H115").Value = 5460.1875
$ws.Range("I115").Value = 3040
$ws.Range("J115").Value = 6560.273
$ws.Range("K115").Value = 9120
$ws.Range("L115").Value = 19680.819
$ws.Range("M115").Value = -7945
$ws.Range("N115").Value = -22030.819
# Row 131
$ws.Range("H131").Value = 14150547
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 14150547
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 42451641
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -42461721

$ws = $wb.Worksheets.Item("GSM")
# Row 140
$ws.Range("H140").Value = 67587.5
$ws.Range("J140").Value = 67587.5
$ws.Range("L140").Value = 67587.5
$ws.Range("N140").Value = -77947.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 771.3158
$ws.Range("I22").Value = 711
$ws.Range("J22").Value = 838.3333
$ws.Range("K22").Value = 711
$ws.Range("L22").Value = 838.3333
$ws.Range("M22").Value = -416
$ws.Range("N22").Value = -1428.3333
# Row 27
$ws.Range("H27").Value = 771.3158
$ws.Range("I27").Value = 711
$ws.Range("J27").Value = 838.3333
$ws.Range("K27").Value = 711
$ws.Range("L27").Value = 838.3333
$ws.Range("M27").Value = -604
$ws.Range("N27").Value = -1052.3333
# Row 46
$ws.Range("H46").Value = 3414.8572
$ws.Range("I46").Value = 2771.1428
$ws.Range("J46").Value = 3575.7856
$ws.Range("K46").Value = 2771.1428
$ws.Range("L46").Value = 3575.7856
$ws.Range("M46").Value = -2583.1428
$ws.Range("N46").Value = -3951.7856
# Row 136
$ws.Range("H136").Value = 3102.4517
$ws.Range("I136").Value = 1996.9333
$ws.Range("J136").Value = 6028.8237
$ws.Range("K136").Value = 5990.7999
$ws.Range("L136").Value = 18086.4711
$ws.Range("M136").Value = -3440.7999
$ws.Range("N136").Value = -23186.4711
